# Edit script: update document text per the target diff.
# Each block finds the original paragraph text and replaces it with the new text,
# preserving the run formatting (bold headings stay bold; body text stays regular).

$d = $word.ActiveDocument

# "2. Dùng AI như vậy có phải là hiệu quả" answer paragraph
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("Không, đây là `"hiệu quả ảo`". Về ngắn hạn, có thể có điểm cao, nhưng sinh viên không học được gì. Khi đi làm thực tế, thiếu nền tảng sẽ khiến bạn ấy bị lộ và mất uy tín. Hiệu quả thực sự là khi AI giúp sinh viên học nhanh hơn, hiểu sâu hơn - chứ không phải thay thế hoàn toàn quá trình tư duy.", $true, $false, $false, $false, $false, $true, 1, $false, "Dùng Ai như vậy sẽ không hiệu quả đối với chúng ta, nó sẽ tạo nên một nghịch lí khá khó nan giải. Về ngắn hạn, có thể có điểm cao, nhưng sinh viên không nạp được những kiến thức từ Ai đưa ra rất dẫn đến việc bản thân người dùng AI sẽ bị thiếu đi kiến thức cơ bản để giải quyết vấn đề nan giải trước mắt. Trong thực tế việc thiếu nền tảng sẽ khiến bản thân rất dễ bị các nhà tuyển dụng nhìn nhận rằng à người phỏng vấn đang thiếu đi kiến thức mà yêu cầu, gây ra làm xấu đi thiện cảm và làm giảm uy tín cho nhà tuyển dụng về bản thân. Hiệu quả thực sự là khi AI giúp sinh viên học nhanh hơn, hiểu sâu hơn - chứ không phải thay thế hoàn toàn quá trình tư duy.", 2)
if (-not $found) { throw "Replacement failed for: para2" }

# "3. Lời khuyên cho sinh viên đó:" heading
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("3. Lời khuyên cho sinh viên đó:", $true, $false, $false, $false, $false, $true, 1, $false, "3. Lời khuyên cho sinh viên đó :", 2)
if (-not $found) { throw "Replacement failed for: para3_heading" }

# "Hãy thành thật..." advice paragraph
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("Hãy thành thật với bản thân về những gì mình chưa biết. Thay vì để AI làm hộ, hãy dùng AI như một người hướng dẫn: hỏi về logic thuật toán, yêu cầu giải thích từng bước, hoặc nhờ review code của mình. Kỹ năng lập trình là tài sản dài hạn, đừng đánh đổi nó lấy điểm số nhất thời.", $true, $false, $false, $false, $false, $true, 1, $false, "Hãy thành thật thói quen với bản thân về những gì mình cảm giác không thể giải thích được. Thay vì sử dụng AI để làm hộ bài tập, hãy dùng AI một cách hợp lí: hỏi vấn đề thuật toán không biết và yêu cầu giải thích từng bước, hoặc nhờ người khác (bạn bè, thầy cô,… ) của mình. Kỹ năng lập trình là tài sản dài hạn, đừng đánh đổi nó lấy điểm số nhất thời.", 2)
if (-not $found) { throw "Replacement failed for: para_hay" }

# "4. Sử dụng AI đúng cách:" heading
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("4. Sử dụng AI đúng cách:", $true, $false, $false, $false, $false, $true, 1, $false, "4. Sử dụng AI một cách hợp lí:", 2)
if (-not $found) { throw "Replacement failed for: para4_heading" }

# "AI nên là công cụ..." paragraph
$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute("AI nên là công cụ hỗ trợ học tập, không phải thay thế. Các cách dùng hợp lý: hỏi AI giải thích khái niệm khó, tham khảo cách tiếp cận bài toán, debug lỗi sau khi đã tự mình cố gắng, hoặc so sánh giải pháp của mình với gợi ý từ AI để rút kinh nghiệm. Luôn tự viết code chính, tự suy nghĩ thuật toán trước khi nhờ AI hỗ trợ", $true, $false, $false, $false, $false, $true, 1, $false, "AI phải là một hỗ trỡ công cụ hỗ trợ học tập, không phải thay thế sức của con người. Những cách dùng hợp lý AI sao cho hợp lí nhất: hỏi AI giải thích khái niệm khó; tham khảo cách tiếp cận bài toán; sửa lỗi sau khi các lỗi được phát sinh ra mà ta cảm giác không làm được việc đó, hoặc so sánh đoạn code của mình với đoạn code đã được AI đề xuất để rút kinh nghiệm. Luôn tự viết code chính, tự suy nghĩ thuật toán trước khi nhờ AI hỗ trợ", 2)
if (-not $found) { throw "Replacement failed for: para_ai_nen" }
